$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text number format on Price cells whose new values would
# otherwise be auto-converted to numeric by Excel (so they stay literal
# text, matching the scraped "Price" column formatting).
$textFormatRows = @(4,5,6,8,9,10,11,13,14,15,17,18,19,20,21,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,43,44,46,47,50,51)
foreach ($r in $textFormatRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) columns with the latest scrape.
$ws.Range("D2").Value = "29.233.58"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "1.860.66"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "0.7028"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "237.73"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.07734"
$ws.Range("E8").Value = "  +4.59%  "
$ws.Range("D9").Value = "0.3049"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").Value = "23.27"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").Value = "0.08191"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("D12").Value = "1.851.12"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "0.7184"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "5.177"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").Value = "89.18"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "29.229.11"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "5.779"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "13.36"
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").Value = "0.000007735"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").Value = "237.29"
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "2.107.43"
$ws.Range("E22").Value = "  +1.55%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "7.435"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("D25").Value = "0.1474"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("D26").Value = "162.09"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "9.011"
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("D29").Value = "2.024"
$ws.Range("E29").Value = "  +5.36%  "
$ws.Range("D30").Value = "1.433"
$ws.Range("E30").Value = "  +4.12%  "
$ws.Range("D31").Value = "4.430"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").Value = "1.487"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "4.056"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("D34").Value = "0.05233"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "1.171"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "0.7083"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "2.668"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").Value = "0.01844"
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("D40").Value = "2.727"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("D41").Value = "0.9324"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("D42").Value = "1.142.43"
$ws.Range("E42").Value = "  +8.00%  "
$ws.Range("D43").Value = "0.4280"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "5.912"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").Value = "103.51"
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("E48").Value = "  +3.08%  "
$ws.Range("D49").Value = "2.004.42"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").Value = "9.177"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").Value = "6.977"
$ws.Range("E51").Value = "  -1.65%  "
